$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# Update existing rows 2-5 with new timestamp/speed/density values
$ws.Range("A2").Value = "2024-08-04 15:11:37"
$ws.Range("B2").Value = 87.94661134592997
$ws.Range("C2").Value = 8

$ws.Range("A3").Value = "2024-08-04 15:11:39"
$ws.Range("B3").Value = 88.85274932208935
$ws.Range("C3").Value = 14

$ws.Range("A4").Value = "2024-08-04 15:11:41"
$ws.Range("B4").Value = 88.54182603881817
$ws.Range("C4").Value = 22

$ws.Range("A5").Value = "2024-08-04 15:11:43"
$ws.Range("B5").Value = 83.86694299403953
$ws.Range("C5").Value = 26

# Append new rows 6-9 with additional simulation data
$ws.Range("A6").Value = "2024-08-04 15:11:45"
$ws.Range("B6").Value = 83.38592155699389
$ws.Range("C6").Value = 34

$ws.Range("A7").Value = "2024-08-04 15:11:47"
$ws.Range("B7").Value = 82.8458131050317
$ws.Range("C7").Value = 36

$ws.Range("A8").Value = "2024-08-04 15:11:49"
$ws.Range("B8").Value = 82.08250678513062
$ws.Range("C8").Value = 37

$ws.Range("A9").Value = "2024-08-04 15:11:51"
$ws.Range("B9").Value = 82.61158489264176
$ws.Range("C9").Value = 41
